$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "시작날짜" (start date) and "종료날짜" (end date) columns (D:E).
# Deleting shifts the old F/G (가중치/학기) columns into D/E, and the old
# L/M columns into J/K, matching the target layout.
$ws.Columns("D:E").Delete() | Out-Null

# Match the author's final cell selection recorded in the saved file.
$ws.Range("G10").Select() | Out-Null
